$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old text classification (A/B/C) in column H to new numeric codes
$map = @{ "A" = 0; "B" = 1; "C     " = 2 }

for ($r = 2; $r -le 65; $r++) {
    $cell = $ws.Cells.Item($r, 8)  # column H
    $text = $cell.Value()
    if ($map.ContainsKey($text)) {
        $cell.Value = $map[$text]
    }
}
